$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INPUT_SHEET")

# ---- Clear the old two data rows, we'll rewrite rows 2-9 from scratch ----
$ws.Range("A2:B3").ClearContents()

# ---- New BAN / ticket data (replaces old row2/row3 content, adds 6 more rows) ----
$data = @(
    @(107198053, "ONREG-20535"),
    @(124473304, "ONREG-18520"),
    @(169068211, "ONREG-19643"),
    @(225356050, "ONREG-22679"),
    @(244834015, "ONREG-22710"),
    @(114478789, "ONREG-25767"),
    @(172500199, "ONREG-21172"),
    @(205653581, "ONREG-26222")
)

$row = 2
foreach ($pair in $data) {
    $ban = $pair[0]
    $ticket = $pair[1]

    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.NumberFormat = "General"
    $aCell.Value = $ban
    $bCell.Value = $ticket

    $aCell.Font.Name = "Times New Roman"
    $aCell.Font.Size = 12
    $aCell.Borders.LineStyle = -4142
    $aCell.HorizontalAlignment = -4152
    $aCell.VerticalAlignment = -4108

    $bCell.Font.Name = "Times New Roman"
    $bCell.Font.Size = 12
    $bCell.Borders.LineStyle = -4142
    $bCell.VerticalAlignment = -4108

    $ws.Rows.Item($row).RowHeight = 15.75

    $row = $row + 1
}

# ---- Trailing formatted blank row (row 10) ----
$a10 = $ws.Cells.Item(10, 1)
$a10.Font.Name = "Segoe UI"
$a10.Font.Size = 10
$a10.Borders.LineStyle = -4142
$a10.VerticalAlignment = -4108

# ---- Column B width ----
$ws.Columns.Item(2).ColumnWidth = 15.43

# ---- Selection matching the author's final cursor position ----
$ws.Range("A9:XFD9").Select()
$ws.Range("B9").Activate()
